# Regenerate orders with updated distance/size codes.
# Distances: D80 -> D86, D51 -> D55, D64 -> D69
# Sizes:     S30 -> S31  (S25, S20 unchanged)
#
# These codes appear embedded inside longer strings (Condition,
# Filename_Left, Filename_Right, Distance, Size columns), so we do a
# substring replace on every used cell's value rather than a whole-value
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = $startRow; $r -lt ($startRow + $rowCount); $r++) {
    for ($c = $startCol; $c -lt ($startCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
